$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp label (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 10:05"

# --- Row 5: Rusia ---
$ws.Range("B5").Value = 326448
$ws.Range("C5").Value = 8894
$ws.Range("D5").Value = 99825
$ws.Range("E5").Value = 223374
$ws.Range("G5").Value = 150
$ws.Range("H5").Value = 3249

# --- Row 14: India ---
$ws.Range("B14").Value = 119419
$ws.Range("C14").Value = 1193
$ws.Range("D14").Value = 48957
$ws.Range("E14").Value = 66863
$ws.Range("G14").Value = 15
$ws.Range("H14").Value = 3599

# --- Rows 90/91: Estonia and Islandia swap ranking position, values updated ---
# Row 90 now holds Estonia's (updated) figures
$ws.Range("A90").Value = "Estonia"
$ws.Range("B90").Value = 1807
$ws.Range("C90").Value = 7
$ws.Range("D90").Value = 1508
$ws.Range("E90").Value = 235
$ws.Range("H90").Value = 64

# Row 91 now holds Islandia's figures
$ws.Range("A91").Value = "Islandia"
$ws.Range("B91").Value = 1803
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 1790
$ws.Range("E91").Value = 3
$ws.Range("H91").Value = 10

# --- Row 97: Eslovaquia ---
$ws.Range("B97").Value = 1503
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 1256
$ws.Range("E97").Value = 219
